# Updated cryptos list values (price + volume change%) and two row swaps
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the numeric-looking price cells as Text so Excel stores them
# as literal strings (matching the source data) instead of coercing them
# into numbers.
$numericTextCells = @("D5","D6","D7","D9","D10","D11","D12","D14","D16","D17","D21","D22","D23","D24","D26","D27","D28","D29","D30","D31","D34","D35","D38","D39","D40","D41","D42","D43","D44","D45","D47","D48","D49")
foreach ($addr in $numericTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "41.855.60"
$ws.Range("E2").Value = "  -4.65%  "

$ws.Range("D3").Value = "2.212.54"
$ws.Range("E3").Value = "  -5.77%  "

$ws.Range("E4").Value = "  -0.42%  "

$ws.Range("D5").Value = "246.28"
$ws.Range("E5").Value = "  +2.72%  "

$ws.Range("D6").Value = "0.628"
$ws.Range("E6").Value = "  -6.06%  "

$ws.Range("D7").Value = "70.81"
$ws.Range("E7").Value = "  -3.67%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").Value = "0.553"
$ws.Range("E9").Value = "  -6.78%  "

$ws.Range("D10").Value = "0.0956"
$ws.Range("E10").Value = "  -5.59%  "

$ws.Range("D11").Value = "58.02"
$ws.Range("E11").Value = "  -5.12%  "

$ws.Range("D12").Value = "36.41"
$ws.Range("E12").Value = "  +8.78%  "

$ws.Range("E13").Value = "  -2.96%  "

$ws.Range("D14").Value = "6.74"
$ws.Range("E14").Value = "  -7.25%  "

$ws.Range("D15").Value = "2.541.67"
$ws.Range("E15").Value = "  -5.83%  "

$ws.Range("D16").Value = "14.91"
$ws.Range("E16").Value = "  -7.75%  "

$ws.Range("D17").Value = "0.846"
$ws.Range("E17").Value = "  -6.35%  "

$ws.Range("D18").Value = "2.210.38"
$ws.Range("E18").Value = "  -5.67%  "

$ws.Range("D19").Value = "41.710.80"
$ws.Range("E19").Value = "  -4.68%  "

$ws.Range("D20").Value = "0.0₃0958"
$ws.Range("E20").Value = "  -6.39%  "

$ws.Range("D21").Value = "73.51"
$ws.Range("E21").Value = "  -5.54%  "

$ws.Range("D22").Value = "6.11"
$ws.Range("E22").Value = "  -7.67%  "

$ws.Range("D23").Value = "235.41"
$ws.Range("E23").Value = "  -6.78%  "

$ws.Range("D24").Value = "2.07"
$ws.Range("E24").Value = "  +12.42%  "

$ws.Range("E25").Value = "  -0.05%  "

$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").Value = "2.46"
$ws.Range("E26").Value = "  -1.26%  "

$ws.Range("B27").Value = "WEMIXToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D27").Value = "3.61"
$ws.Range("E27").Value = "  -5.47%  "

$ws.Range("D28").Value = "2.24"
$ws.Range("E28").Value = "  -0.76%  "

$ws.Range("D29").Value = "9.93"
$ws.Range("E29").Value = "  -4.57%  "

$ws.Range("D30").Value = "169.87"
$ws.Range("E30").Value = "  -3.41%  "

$ws.Range("D31").Value = "20.47"
$ws.Range("E31").Value = "  -7.81%  "

$ws.Range("E32").Value = "  -5.89%  "

$ws.Range("E33").Value = "  -7.05%  "

$ws.Range("D34").Value = "0.0717"
$ws.Range("E34").Value = "  -3.49%  "

$ws.Range("D35").Value = "5.15"
$ws.Range("E35").Value = "  -3.85%  "

$ws.Range("E37").Value = "  +2.07%  "

$ws.Range("D38").Value = "23.30"
$ws.Range("E38").Value = "  +20.35%  "

$ws.Range("D39").Value = "2.28"
$ws.Range("E39").Value = "  -4.92%  "

$ws.Range("D40").Value = "0.0273"
$ws.Range("E40").Value = "  +0.26%  "

$ws.Range("D41").Value = "5.88"
$ws.Range("E41").Value = "  -8.41%  "

$ws.Range("D42").Value = "66.22"
$ws.Range("E42").Value = "  +0.71%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "8.98"
$ws.Range("E43").Value = "  -1.76%  "

$ws.Range("B44").Value = "FTXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D44").Value = "4.85"
$ws.Range("E44").Value = "  -9.82%  "

$ws.Range("B45").Value = "SynthetixNetwork"
$ws.Range("C45").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D45").Value = "4.64"
$ws.Range("E45").Value = "  +7.98%  "

$ws.Range("E46").Value = "  -5.30%  "

$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "0.0996"
$ws.Range("E47").Value = "  -5.32%  "

$ws.Range("B48").Value = "BinanceUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D48").Value = "1.00"
$ws.Range("E48").Value = "  -0.04%  "

$ws.Range("D49").Value = "10.34"
$ws.Range("E49").Value = "  +8.22%  "

$ws.Range("E50").Value = "  -3.51%  "

$ws.Range("E51").Value = "  -4.11%  "
